# schedule_aller.xlsx — "add colors TOPAZE and AMETHYSTE"
#
# - Header band (row 1) recolored from BLUE (4472C4) to AMETHYSTE (9966CC),
#   and its font gains an explicit Arial face (it was already bold/white).
# - The two alternating zebra-stripe bands used by the data rows move from
#   the old blue family to a purple family:
#     light-blue BDD7EE -> lavender  E6E6FA
#     light-green E2EFDA -> thistle  D8BFD8
# - Column F narrows slightly (10 -> 9 chars).
# - Header captions are translated/relabelled to French, capitalised.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- color constants (OLE COLORREF = R + G*256 + B*65536) ----
$COLOR_AMETHYSTE         = 13395609  # 9966CC
$COLOR_OLD_LIGHTBLUE     = 15652797  # BDD7EE
$COLOR_LAVENDER          = 16443110  # E6E6FA
$COLOR_OLD_LIGHTGREEN    = 14348258  # E2EFDA
$COLOR_THISTLE           = 14204888  # D8BFD8

# ---- header row: relabel + restyle ----
$ws.Range("A1").Value = "Round"
$ws.Range("B1").Value = "Début"
$ws.Range("C1").Value = "Fin"
$ws.Range("D1").Value = "Équipe 1"
$ws.Range("E1").Value = "Équipe 2"
$ws.Range("F1").Value = "Durée"
$ws.Range("G1").Value = "Phase"

$header = $ws.Range("A1:G1")
$header.Font.Name = "Arial"
$header.Interior.Color = $COLOR_AMETHYSTE

# ---- data rows: remap the two zebra-stripe bands to the purple family ----
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $row = $ws.Range("A" + $r + ":G" + $r)
    $existing = $ws.Cells.Item($r, 1).Interior.Color

    if ($existing -eq $COLOR_OLD_LIGHTBLUE) {
        $row.Interior.Color = $COLOR_LAVENDER
    } elseif ($existing -eq $COLOR_OLD_LIGHTGREEN) {
        $row.Interior.Color = $COLOR_THISTLE
    }
}

# ---- column F narrows from 10 to 9 characters ----
$ws.Columns.Item(6).ColumnWidth = 9
